$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 8959.214
$ws.Range("J17").Value = 9533
$ws.Range("L17").Value = 28599
$ws.Range("N17").Value = -28935
$ws.Range("H42").Value = 524.25
$ws.Range("I42").Value = 370
$ws.Range("K42").Value = 1110
$ws.Range("M42").Value = -880
$ws.Range("H46").Value = 3159.6
$ws.Range("I46").Value = 3266
$ws.Range("J46").Value = 3000
$ws.Range("K46").Value = 9798
$ws.Range("L46").Value = 9000
$ws.Range("M46").Value = -9679
$ws.Range("N46").Value = -9238
$ws.Range("H60").Value = 3159.6
$ws.Range("I60").Value = 3266
$ws.Range("J60").Value = 3000
$ws.Range("K60").Value = 9798
$ws.Range("L60").Value = 9000
$ws.Range("M60").Value = -9314
$ws.Range("N60").Value = -9968
$ws.Range("H112").Value = 2060.8635
$ws.Range("J112").Value = 2087.5715
$ws.Range("L112").Value = 6262.7145
$ws.Range("N112").Value = -8478.7145
$ws.Range("H116").Value = 5353.2666
$ws.Range("I116").Value = 3649.5
$ws.Range("J116").Value = 7300.4287
$ws.Range("K116").Value = 3649.5
$ws.Range("L116").Value = 7300.4287
$ws.Range("M116").Value = -207.5
$ws.Range("N116").Value = -14184.4287
$ws.Range("H132").Value = 960.8837
$ws.Range("I132").Value = 987.975
$ws.Range("K132").Value = 2963.925
$ws.Range("M132").Value = -433.9250000000002

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 3008
$ws.Range("I22").Value = 3008
$ws.Range("K22").Value = 3008
$ws.Range("M22").Value = -2709
$ws.Range("H32").Value = 4550.7935
$ws.Range("I32").Value = 3251.5862
$ws.Range("J32").Value = 19621.6
$ws.Range("K32").Value = 3251.5862
$ws.Range("L32").Value = 19621.6
$ws.Range("M32").Value = -2964.5862
$ws.Range("N32").Value = -20195.6
$ws.Range("H88").Value = 4108.684
$ws.Range("I88").Value = 3280.8333
$ws.Range("J88").Value = 5527.857
$ws.Range("K88").Value = 3280.8333
$ws.Range("L88").Value = 5527.857
$ws.Range("M88").Value = -2874.8333
$ws.Range("N88").Value = -6339.857
$ws.Range("H91").Value = 4108.684
$ws.Range("I91").Value = 3280.8333
$ws.Range("J91").Value = 5527.857
$ws.Range("K91").Value = 3280.8333
$ws.Range("L91").Value = 5527.857
$ws.Range("M91").Value = -1876.8333
$ws.Range("N91").Value = -8335.857
$ws.Range("H132").Value = 5205.4287
$ws.Range("I132").Value = 3043.6667
$ws.Range("K132").Value = 9131.000100000001
$ws.Range("M132").Value = -6601.000100000001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4242.303
$ws.Range("I20").Value = 3960.4783
$ws.Range("K20").Value = 3960.4783
$ws.Range("M20").Value = -3713.4783
$ws.Range("H37").Value = 6157.6
$ws.Range("I37").Value = 6508.4443
$ws.Range("K37").Value = 6508.4443
$ws.Range("M37").Value = -6371.4443
$ws.Range("H94").Value = 833.2273
$ws.Range("I94").Value = 833.2273
$ws.Range("K94").Value = 833.2273
$ws.Range("M94").Value = -382.2273
$ws.Range("H96").Value = 20000
$ws.Range("I96").Value = 20000
$ws.Range("K96").Value = 20000
$ws.Range("M96").Value = -17254
$ws.Range("H134").Value = 2197.0244
$ws.Range("I134").Value = 1411.1111
$ws.Range("K134").Value = 4233.3333
$ws.Range("M134").Value = -1698.3333

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 59578
$ws.Range("J20").Value = 59578
$ws.Range("L20").Value = 59578
$ws.Range("N20").Value = -60050
$ws.Range("H29").Value = 4921
$ws.Range("J29").Value = 4921
$ws.Range("L29").Value = 4921
$ws.Range("N29").Value = -5507
$ws.Range("H30").Value = 59578
$ws.Range("J30").Value = 59578
$ws.Range("L30").Value = 59578
$ws.Range("N30").Value = -59760
$ws.Range("H31").Value = 43620.27
$ws.Range("I31").Value = 3154.1667
$ws.Range("K31").Value = 3154.1667
$ws.Range("M31").Value = -2859.1667
$ws.Range("H34").Value = 43620.27
$ws.Range("I34").Value = 3154.1667
$ws.Range("K34").Value = 3154.1667
$ws.Range("M34").Value = -2952.1667
$ws.Range("H122").Value = 4039.5417
$ws.Range("I122").Value = 1439.3572
$ws.Range("K122").Value = 4318.071599999999
$ws.Range("M122").Value = -1868.071599999999
$ws.Range("H128").Value = 59578
$ws.Range("J128").Value = 59578
$ws.Range("L128").Value = 59578
$ws.Range("N128").Value = -69538
$ws.Range("H132").Value = 2218.8806
$ws.Range("I132").Value = 1827.5167
$ws.Range("K132").Value = 5482.550099999999
$ws.Range("M132").Value = -2952.550099999999

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 1321.3334
$ws.Range("I45").Value = 1032
$ws.Range("J45").Value = 1466
$ws.Range("K45").Value = 3096
$ws.Range("L45").Value = 4398
$ws.Range("M45").Value = -2564
$ws.Range("N45").Value = -5462
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H109").Value = 2619.0715
$ws.Range("I109").Value = 1265.625
$ws.Range("J109").Value = 4423.6665
$ws.Range("K109").Value = 3796.875
$ws.Range("L109").Value = 13270.9995
$ws.Range("M109").Value = -2756.875
$ws.Range("N109").Value = -15350.9995
$ws.Range("H131").Value = 9984423
$ws.Range("I131").Value = 41667380
$ws.Range("J131").Value = 6706875
$ws.Range("K131").Value = 125002140
$ws.Range("L131").Value = 20120625
$ws.Range("M131").Value = -124997100
$ws.Range("N131").Value = -20130705

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10105.071
$ws.Range("I70").Value = 9416.166999999999
$ws.Range("J70").Value = 10621.75
$ws.Range("K70").Value = 9416.166999999999
$ws.Range("L70").Value = 10621.75
$ws.Range("M70").Value = -9146.166999999999
$ws.Range("N70").Value = -11161.75
$ws.Range("H73").Value = 10105.071
$ws.Range("I73").Value = 9416.166999999999
$ws.Range("J73").Value = 10621.75
$ws.Range("K73").Value = 9416.166999999999
$ws.Range("L73").Value = 10621.75
$ws.Range("M73").Value = -8480.166999999999
$ws.Range("N73").Value = -12493.75
$ws.Range("H113").Value = 3052.6
$ws.Range("I113").Value = 2191.9167
$ws.Range("K113").Value = 2191.9167
$ws.Range("M113").Value = -21.91670000000022
$ws.Range("H136").Value = 19581.959
$ws.Range("J136").Value = 19581.959
$ws.Range("L136").Value = 58745.87699999999
$ws.Range("N136").Value = -63845.87699999999

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6065121.5
$ws.Range("I7").Value = 8698983
$ws.Range("J7").Value = 7240.6
$ws.Range("K7").Value = 8698983
$ws.Range("L7").Value = 7240.6
$ws.Range("M7").Value = -8698871
$ws.Range("N7").Value = -7464.6
$ws.Range("H22").Value = 2591.4644
$ws.Range("I22").Value = 1130.8334
$ws.Range("J22").Value = 3686.9375
$ws.Range("K22").Value = 1130.8334
$ws.Range("L22").Value = 3686.9375
$ws.Range("M22").Value = -835.8334
$ws.Range("N22").Value = -4276.9375
$ws.Range("H27").Value = 2591.4644
$ws.Range("I27").Value = 1130.8334
$ws.Range("J27").Value = 3686.9375
$ws.Range("K27").Value = 1130.8334
$ws.Range("L27").Value = 3686.9375
$ws.Range("M27").Value = -1023.8334
$ws.Range("N27").Value = -3900.9375
$ws.Range("H93").Value = 3736.3333
$ws.Range("I93").Value = 3745.6
$ws.Range("K93").Value = 3745.6
$ws.Range("M93").Value = -2497.6
$ws.Range("H122").Value = 152798.97
$ws.Range("I122").Value = 203998.5
$ws.Range("K122").Value = 611995.5
$ws.Range("M122").Value = -609545.5
$ws.Range("H126").Value = 6065121.5
$ws.Range("I126").Value = 8698983
$ws.Range("J126").Value = 7240.6
$ws.Range("K126").Value = 26096949
$ws.Range("L126").Value = 21721.8
$ws.Range("M126").Value = -26094479
$ws.Range("N126").Value = -26661.8

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 24875
$ws.Range("I54").Value = 24500
$ws.Range("J54").Value = 25000
$ws.Range("K54").Value = 24500
$ws.Range("L54").Value = 25000
$ws.Range("M54").Value = -23980
$ws.Range("N54").Value = -26040
$ws.Range("H81").Value = 3683.7058
$ws.Range("J81").Value = 4769.5
$ws.Range("L81").Value = 9539
$ws.Range("N81").Value = -11661
$ws.Range("H84").Value = 3683.7058
$ws.Range("J84").Value = 4769.5
$ws.Range("L84").Value = 47695
$ws.Range("N84").Value = -58303
$ws.Range("H136").Value = 2213.2273
$ws.Range("I136").Value = 499.13794
$ws.Range("K136").Value = 1497.41382
$ws.Range("M136").Value = 1052.58618

Write-Host "Applied 218 cell updates and 1 clears"